# The course schedule listed several sections (AV-..., AM-..., AT-... prefixed
# courses) that no longer belong on the sheet - likely duplicated rows that
# leaked in because a pandas DataFrame was reused by reference instead of
# being copied ("Dataframes pass by reference when you create a DF from
# another DF"). This script strips those stray "AV-"/"AM-"/"AT-" course
# lines out of every multi-line cell on the single worksheet, leaving the
# period number and any remaining CM-/CSM- course lines intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $val -is [string] -and $val.Contains("`n")) {
            $lines = $val -split "`n"
            $keep = @($lines[0])
            for ($i = 1; $i -lt $lines.Length; $i++) {
                $line = $lines[$i]
                if (-not ($line.StartsWith("AV-") -or $line.StartsWith("AM-") -or $line.StartsWith("AT-"))) {
                    $keep = $keep + @($line)
                }
            }
            $newVal = $keep -join "`n"
            if ($newVal -ne $val) {
                if ($keep.Length -eq 1 -and $newVal -match '^-?\d+(\.\d+)?$') {
                    # The remaining text is just the bare period number (e.g. "13",
                    # "8.0"...) with no course line left beneath it. Left as a plain
                    # assignment Excel would silently re-type the cell as a number;
                    # force it to stay text (matching the original shared-string
                    # cell type) the same way typing into a pre-formatted Text cell
                    # would.
                    $cell.NumberFormat = "@"
                }
                $cell.Value = $newVal
            }
        }
    }
}
